$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.222.39'
$ws.Range("E2").Value = '  +5.38%  '
$ws.Range("D3").Value = '2.240.11'
$ws.Range("E3").Value = '  +3.11%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.44'
$ws.Range("E5").Value = '  +4.26%  '
$ws.Range("E6").Value = '  +1.38%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '76.28'
$ws.Range("E7").Value = '  +9.47%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.616'
$ws.Range("E9").Value = '  +7.43%  '
$ws.Range("E10").Value = '  +4.25%  '
$ws.Range("E11").Value = '  +1.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.83'
$ws.Range("E12").Value = '  +1.32%  '
$ws.Range("E13").Value = '  +3.69%  '
$ws.Range("E14").Value = '  +0.77%  '
$ws.Range("D15").Value = '2.572.86'
$ws.Range("E15").Value = '  +3.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.84'
$ws.Range("E16").Value = '  +7.34%  '
$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.817'
$ws.Range("E17").Value = '  +0.97%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '2.233.69'
$ws.Range("E18").Value = '  +3.09%  '
$ws.Range("D19").Value = '43.105.19'
$ws.Range("E19").Value = '  +5.53%  '
$ws.Range("E20").Value = '  +3.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.15'
$ws.Range("E21").Value = '  +1.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.02'
$ws.Range("E22").Value = '  +1.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.61'
$ws.Range("E23").Value = '  +8.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '231.08'
$ws.Range("E24").Value = '  +2.60%  '
$ws.Range("E25").Value = '  +10.79%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.98'
$ws.Range("E27").Value = '  +0.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.35'
$ws.Range("E28").Value = '  -5.40%  '
$ws.Range("E29").Value = '  +2.59%  '
$ws.Range("E30").Value = '  +0.41%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '174.54'
$ws.Range("E31").Value = '  +5.22%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '37.41'
$ws.Range("E32").Value = '  +21.06%  '
$ws.Range("E33").Value = '  +2.84%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0797'
$ws.Range("E34").Value = '  +3.06%  '
$ws.Range("E35").Value = '  +5.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.123'
$ws.Range("E36").Value = '  +1.62%  '
$ws.Range("E37").Value = '  +7.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.34'
$ws.Range("E38").Value = '  +4.74%  '
$ws.Range("E39").Value = '  +15.50%  '
$ws.Range("E40").Value = '  +6.66%  '
$ws.Range("E41").Value = '  +4.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.61'
$ws.Range("E42").Value = '  +3.32%  '
$ws.Range("E43").Value = '  +5.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '60.46'
$ws.Range("E44").Value = '  +1.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '105.60'
$ws.Range("E45").Value = '  +6.77%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.55'
$ws.Range("E46").Value = '  +3.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0996'
$ws.Range("E47").Value = '  +2.48%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.457'
$ws.Range("E48").Value = '  +24.64%  '
$ws.Range("E49").Value = '  +2.66%  '
$ws.Range("E50").Value = '  +4.46%  '
$ws.Range("E51").Value = '  +2.25%  '
